# Apply the NP_Database.xlsx edit:
#  - Sheet "NPUNID": drop the trailing placeholder rows (52-101) that only
#    carried an index number in column A.
#  - Sheet "NP_Props" / Table1: rename the "Ligand" column to "Ligands",
#    widen it, and replace several of its values with the fuller
#    "Core_Shell..." style ligand description used elsewhere in the sheet;
#    also backfill the newly-visible "Dh_core" (E) values for the rows that
#    previously had an empty cell there.
#  - Make "NP_Props" the active/selected sheet, restoring the previous
#    selections on both sheets.

$wb = $excel.ActiveWorkbook

$wsNPUNID = $wb.Worksheets.Item("NPUNID")
$wsProps  = $wb.Worksheets.Item("NP_Props")

# ---------------------------------------------------------------------
# Sheet "NPUNID": remove the now-unused tail rows (52-101), which only
# held sequential index numbers in column A and nothing else.
# ---------------------------------------------------------------------
$wsNPUNID.Range("A52:A101").EntireRow.Delete()

# ---------------------------------------------------------------------
# Sheet "NP_Props" (Table1): rename column D header, which also renames
# the table column automatically.
# ---------------------------------------------------------------------
$wsProps.Range("D1").Value = "Ligands"

# Widen column D to fit the longer ligand descriptions.
$wsProps.Columns.Item(4).ColumnWidth = 91.2

# Update the ligand descriptions (column D) for the gold-core rows and
# fill in the companion Dh_core (column E) values that were previously
# blank.
$wsProps.Range("E11").Value = 200
$wsProps.Range("E12").Value = 200

$wsProps.Range("E13").Value = 71.5
$wsProps.Range("E14").Value = 200

$wsProps.Range("D15").Value = "Citrate_Polyethelyenimine"
$wsProps.Range("E15").Value = 71.5
$wsProps.Range("D16").Value = "Citrate_Polyethelyenimine"
$wsProps.Range("E16").Value = 200

$wsProps.Range("D17").Value = "Citrate_Polyethelyenimine_Au_polyvinylpyrrolidone"
$wsProps.Range("E17").Value = 200
$wsProps.Range("D18").Value = "Citrate_Polyethelyenimine_Au_polyvinylpyrrolidone"
$wsProps.Range("E18").Value = 71.5
$wsProps.Range("D19").Value = "Citrate_Polyethelyenimine_Au_polyvinylpyrrolidone"
$wsProps.Range("E19").Value = 71.5

$wsProps.Range("E20").Value = 200

$wsProps.Range("D21").Value = "Citrate_Polyethelyenimine"
$wsProps.Range("E21").Value = 200

$wsProps.Range("E22").Value = 71.5

$wsProps.Range("D23").Value = "Citrate_Polyethelyenimine"
$wsProps.Range("E23").Value = 71.5

$wsProps.Range("D24").Value = "Citrate_Polyethelyenimine_Au_Polyvinylpyrrolidone"
$wsProps.Range("E24").Value = 200

$wsProps.Range("D25").Value = "Citrate_Polyethelyenimine_Au_PEG5k"
$wsProps.Range("E25").Value = 200

$wsProps.Range("D26").Value = "Citrate_Polyethelyenimine_Au_Polyethelyenimine"
$wsProps.Range("E26").Value = 200
$wsProps.Range("D27").Value = "Citrate_Polyethelyenimine_Au_Polyethelyenimine"
$wsProps.Range("E27").Value = 71.5

$wsProps.Range("D28").Value = "Carboxylate_PEG2k"
$wsProps.Range("E28").Value = 200

# ---------------------------------------------------------------------
# Selections / active sheet: NP_Props becomes the active (selected) tab,
# with its own selection moved to D23; NPUNID keeps a selection at D16.
# ---------------------------------------------------------------------
$wsNPUNID.Range("D16").Select()
$wsProps.Activate()
$wsProps.Range("D23").Select()
